# Insert a new weekly price record for "Flame Seedless" grapes at Terminal
# La Palmera de La Serena. This pushes the existing rows 119-130 down to
# 120-131 (dimension grows from A1:T130 to A1:T131).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 119, shifting the rest down.
$ws.Rows(119).Insert()

# Populate the newly inserted row 119 with the new record's data.
$ws.Range("A119").Value = 8
$ws.Range("B119").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C119").Value = 'Coquimbo'
$ws.Range("D119").Value = 44931
$ws.Range("E119").Value = 4
$ws.Range("F119").Value = 'Fruta'
$ws.Range("G119").Value = 100109
$ws.Range("H119").Value = 'Uva'
$ws.Range("I119").Value = 100109001
$ws.Range("J119").Value = 'Uva'
$ws.Range("K119").Value = 'Flame Seedless'
$ws.Range("L119").Value = 'Primera'
$ws.Range("M119").Value = 800
$ws.Range("N119").Value = 9000
$ws.Range("O119").Value = 10000
$ws.Range("P119").Value = 9500
$ws.Range("Q119").Value = '$/bandeja 10 kilos'
$ws.Range("R119").Value = 'Provincia de Limarí'
$ws.Range("S119").Value = 950
$ws.Range("T119").Value = 10
